# Applies the "Student Book a class General" edit to the Notes workbook.
#
# Summary of changes (see commit diff):
#  - Row 10/11: role column (B) "tutor" -> "d"
#  - Row 12: role column (B) "d" -> "tutor"
#  - Row 10/11/12: username column (D) -> "tutor60@nkt.com"
#  - Row 12: password column (E) -> "Admin@123"
#  - Row 10/11/12: class_name column (F) -> "Pilot sess 6" / "Pilot MV 6" / "Pilot MF 6"
#  - Hyperlinks on E10:E12 merged into a single hyperlink (mailto:Admin@123)
#  - New hyperlink on D10 (mailto:tutor60@nkt.com)
#  - Hyperlinks on D11:D12 merged into a single hyperlink (mailto:tutor60@nkt.com)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# ---------------------------------------------------------------------------
# Helper: delete the hyperlink anchored at a specific single-cell address.
# (Range.Hyperlinks.Delete() removes every hyperlink on the sheet in this
# runtime, so we locate + delete the individual Hyperlink object instead.)
# ---------------------------------------------------------------------------
function Remove-HyperlinkAt($sheet, $address) {
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Address() -eq $address) {
            $h.Delete()
            return $true
        }
    }
    return $false
}

# Remove the three existing (soon to be merged/replaced) hyperlinks.
Remove-HyperlinkAt $ws '$E$12' | Out-Null
Remove-HyperlinkAt $ws '$E$11' | Out-Null
Remove-HyperlinkAt $ws '$E$10' | Out-Null

# ---------------------------------------------------------------------------
# Update cell values for rows 10-12.
# ---------------------------------------------------------------------------

# Row 10
$ws.Range("B10").Value = "d"
$ws.Range("D10").Value = "tutor60@nkt.com"
$ws.Range("F10").Value = "Pilot sess 6"

# Row 11
$ws.Range("B11").Value = "d"
$ws.Range("D11").Value = "tutor60@nkt.com"
$ws.Range("F11").Value = "Pilot MV 6"

# Row 12
$ws.Range("B12").Value = "tutor"
$ws.Range("D12").Value = "tutor60@nkt.com"
$ws.Range("E12").Value = "Admin@123"
$ws.Range("F12").Value = "Pilot MF 6"

# ---------------------------------------------------------------------------
# Re-create the hyperlinks with their new ranges/targets.
# ---------------------------------------------------------------------------

# E10:E12 -> mailto:Admin@123 (single hyperlink spanning the merged range)
$ws.Hyperlinks.Add($ws.Range("E10:E12"), "mailto:Admin@123", "", "", "Admin@123") | Out-Null

# D10 -> mailto:tutor60@nkt.com (single cell, text already matches so no display override)
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:tutor60@nkt.com") | Out-Null

# D11:D12 -> mailto:tutor60@nkt.com (single hyperlink spanning the merged range)
$ws.Hyperlinks.Add($ws.Range("D11:D12"), "mailto:tutor60@nkt.com", "", "", "tutor60@nkt.com") | Out-Null

# Give D10/D11/D12 the same "Hyperlink" cell style already used by D7:D9 so the
# username column keeps a consistent look.
$ws.Range("D10:D12").Style = "Hyperlink"
